$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaktiomaksut")

# kayttotarkoitus-hinnasto paivitys:
# 711 teollisuusvarastot          -> A
# 699 muut teollisuuden tuotantorakennukset -> D
# 613 yhdyskuntatekniikan rakennukset       -> D
$ws.Range("B59").Value = "D"
$ws.Range("B62").Value = "D"
$ws.Range("B63").Value = "A"
